$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6921212077140808
$ws.Range("B1").Value = 0.9709581136703491
$ws.Range("C1").Value = 1.282326340675354
$ws.Range("D1").Value = 4.12547492980957
$ws.Range("E1").Value = 2.416555643081665
